$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.176.27"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.76"
$ws.Range("E3").Value = "  +3.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.75"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4690"
$ws.Range("E7").Value = "  +24.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3715"
$ws.Range("E8").Value = "  +11.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.30"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07709"
$ws.Range("E10").Value = "  +7.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.158"
$ws.Range("E11").Value = "  +4.46%  "

$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.389"
$ws.Range("E14").Value = "  +4.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.407"
$ws.Range("E15").Value = "  +4.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.799.26"
$ws.Range("E16").Value = "  +3.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  +4.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06736"
$ws.Range("E18").Value = "  +2.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.74"
$ws.Range("E19").Value = "  +4.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.51"
$ws.Range("E21").Value = "  +4.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.444"
$ws.Range("E22").Value = "  +3.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.159.27"
$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.95"
$ws.Range("E24").Value = "  +3.02%  "

$ws.Range("E25").Value = "  +0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.99"
$ws.Range("E26").Value = "  +6.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.410"
$ws.Range("E27").Value = "  +4.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.22"
$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.007.44"
$ws.Range("E29").Value = "  +3.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.65"
$ws.Range("E30").Value = "  +2.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.276"
$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.045"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09658"
$ws.Range("E33").Value = "  +10.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.943"
$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2246"
$ws.Range("E35").Value = "  +6.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.30"
$ws.Range("E36").Value = "  +1.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02384"
$ws.Range("E37").Value = "  +3.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06419"
$ws.Range("E38").Value = "  +3.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6756"
$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.275"
$ws.Range("E40").Value = "  +2.50%  "

$ws.Range("E41").Value = "  +5.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.238"
$ws.Range("E42").Value = "  +2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.153"
$ws.Range("E43").Value = "  +2.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.17"
$ws.Range("E44").Value = "  +3.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6204"
$ws.Range("E46").Value = "  +3.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.839"
$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.18"
$ws.Range("E48").Value = "  +2.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.075"
$ws.Range("E49").Value = "  +3.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.191"
$ws.Range("E50").Value = "  +3.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07155"
$ws.Range("E51").Value = "  +0.81%  "
